$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 85; this shifts rows 85:152 down to 86:153
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new record's data
$ws.Cells.Item(85, 1).Value = 7
$ws.Cells.Item(85, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(85, 3).Value = "Ñuble"
$ws.Cells.Item(85, 4).Value = 44447
$ws.Cells.Item(85, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(85, 5).Value = 16
$ws.Cells.Item(85, 6).Value = 100112023
$ws.Cells.Item(85, 7).Value = "Brócoli"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 300
$ws.Cells.Item(85, 11).Value = 750
$ws.Cells.Item(85, 12).Value = 800
$ws.Cells.Item(85, 13).Value = 775
$ws.Cells.Item(85, 14).Value = "$/unidad"
$ws.Cells.Item(85, 15).Value = "Región del Maule"
$ws.Cells.Item(85, 16).Value = 775
$ws.Cells.Item(85, 17).Value = 1
$ws.Cells.Item(85, 18).Value = "Hortaliza"
